$d = $word.ActiveDocument

# The document has a single paragraph whose whole text is "luis" in one run.
# Target: split that text into three runs -> "L" | "uis" | "-s-sdsdsdss"
# (the last run appends the literal "-s-sdsdsdss" suffix).

$para = $d.Paragraphs.First
$pr = $para.Range

# Range covering only the paragraph's content, excluding the trailing
# paragraph mark, so paragraph-level attributes (paraId/textId/rsid...)
# stay untouched.
$contentRange = $d.Range($pr.Start, $pr.End - 1)

$xml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
       "<w:r><w:t>L</w:t></w:r>" +
       "<w:r><w:t>uis</w:t></w:r>" +
       "<w:r><w:t>-s-sdsdsdss</w:t></w:r>" +
       "</w:p>"

$contentRange.InsertXML($xml)
